$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Row 6
$ws.Range("D6").Value = ""
$ws.Range("H6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("I6").Value = ""
$ws.Range("K6").Value = ""

# Row 12
$ws.Range("D12").Value = ""
$ws.Range("K12").Value = ""

# Row 13
$ws.Range("D13").Value = ""
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 7,97 TL"
$ws.Range("I13").Value = ""
$ws.Range("K13").Value = ""

# Row 14
$ws.Range("D14").Value = ""
$ws.Range("H14").Value = "3.000 TL - 6.000 TL"
$ws.Range("K14").Value = ""
